$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold numeric-/percentage-looking values
# that are stored as plain text in the workbook. Excel would otherwise parse
# strings such as "285.98" or "2.46%" as real numbers (and apply a General /
# Percentage number format), so each such cell is temporarily switched to the
# Text ("@") number format before the value is written, then the formatting is
# cleared again so the cell keeps its original (unstyled) appearance while the
# value remains literal text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"

$ws.Range("D2").Value = '285.98'
$ws.Range("E2").Value = '2.46%'
$ws.Range("D3").Value = '28.81'
$ws.Range("E3").Value = '5.06%'
$ws.Range("D4").Value = '5.046'
$ws.Range("E4").Value = '4.28%'
$ws.Range("D5").Value = '0.06702'
$ws.Range("E5").Value = '5.05%'
$ws.Range("D6").Value = '7.343'
$ws.Range("E6").Value = '4.38%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '3.383'
$ws.Range("E7").Value = '2.33%'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = '1.377'
$ws.Range("E8").Value = '6.79%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '0.9390'
$ws.Range("E9").Value = '5.05%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1567'
$ws.Range("E10").Value = '2.86%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '0.06674'
$ws.Range("E11").Value = '14.29%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.07577'
$ws.Range("E12").Value = '0.88%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.02964'
$ws.Range("E13").Value = '1.00%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.08984'
$ws.Range("E14").Value = '-0.08%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001587'
$ws.Range("E15").Value = '1.34%'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = '0.04491'
$ws.Range("E16").Value = '1.96%'
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").Value = '0.0006447'
$ws.Range("E17").Value = '0.73%'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Value = '0.006502'
$ws.Range("E18").Value = '7.07%'
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").Value = '3.493'
$ws.Range("E19").Value = '0.44%'
$ws.Range("E20").Value = '0.90%'
$ws.Range("D21").Value = '0.3206'
$ws.Range("E21").Value = '1.89%'
$ws.Range("E22").Value = '-3.04%'
$ws.Range("D23").Value = '4.088'
$ws.Range("E23").Value = '4.73%'
$ws.Range("D24").Value = '0.1550'
$ws.Range("E24").Value = '3.16%'
$ws.Range("D25").Value = '0.001181'
$ws.Range("E25").Value = '0.35%'
$ws.Range("E26").Value = '4.99%'
$ws.Range("D27").Value = '0.0001247'
$ws.Range("E27").Value = '5.67%'
$ws.Range("E28").Value = '-2.28%'
$ws.Range("D40").Value = '0.04214'
$ws.Range("E40").Value = '3.29%'
$ws.Range("D41").Value = '0.006731'
$ws.Range("E41").Value = '1.97%'
$ws.Range("D42").Value = '0.1252'
$ws.Range("E42").Value = '-11.05%'
$ws.Range("D43").Value = '0.002015'
$ws.Range("E43").Value = '-5.40%'
$ws.Range("D44").Value = '0.01222'
$ws.Range("E44").Value = '11.12%'
$ws.Range("D45").Value = '0.00005633'
$ws.Range("E45").Value = '1.92%'
$ws.Range("D47").Value = '0.01304'
$ws.Range("E47").Value = '-29.42%'

$ws.Range("D2").ClearFormats()
$ws.Range("E2").ClearFormats()
$ws.Range("D3").ClearFormats()
$ws.Range("E3").ClearFormats()
$ws.Range("D4").ClearFormats()
$ws.Range("E4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("E5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("E6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("E7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("E8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("E9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("E10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("E11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("E12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("E13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("E14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("E15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("E16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("E17").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("E18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("E19").ClearFormats()
$ws.Range("E20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("E21").ClearFormats()
$ws.Range("E22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("E23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("E24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("E25").ClearFormats()
$ws.Range("E26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("E27").ClearFormats()
$ws.Range("E28").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("E40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("E41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("E42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("E43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("E44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("E45").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("E47").ClearFormats()
